$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows to open up the new layout:
#  - insert a blank row at 40 (pushes old 40 "ECONOMIC CENSUS..." text down to 41,
#    old hyperlink row 41 down to 42, etc.)
$ws.Rows("40").Insert()
#  - remove the row that is now a duplicate blank separator at 42
$ws.Rows("42").Delete()
#  - insert a new row at 43 for the (former hyperlink) URL text, now plain text
$ws.Rows("43").Insert()

# Row 41: was the hyperlinked URL text; now just plain "source"-styled text
$ws.Range("A41").Value = "ECONOMIC CENSUS 2011 - PROFILE OF SMALL AND MEDIUM ENTERPRISE, Department of Statistics Malaysia"

# Row 43: the URL text, still styled like the other source lines (no hyperlink)
$ws.Range("A43").Value = "http://www.statistics.gov.my/portal/index.php?option=com_content&view=article&id=1721&Itemid=149&lang=en"

# Remove the now-obsolete hyperlink object (text stays, formatting no longer blue/underlined)
foreach ($h in $ws.Hyperlinks) {
    $h.Delete()
}

# Row 47 (previously the DOSM economic-census citation) now references the new SME definition guideline
$ws.Range("A47").Value = "GUIDELINE FOR NEW SME DEFINITION, 2013 available at http://www.smecorp.gov.my/vn2/sites/default/files/Guideline_for_New_SME_Definition_7Jan2014.pdf"

# Row 49 (previously the informal-sector survey citation) now just reads "DOSM"
$ws.Range("A49").Value = "DOSM"
